$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 209 (shifts existing rows 209..291 down to 210..292)
$ws.Rows.Item(209).Insert()

# Populate the newly inserted row 209 with the new weekly record
$ws.Cells.Item(209, 1).Value = 11
$ws.Cells.Item(209, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(209, 3).Value = "Bíobío"
$ws.Cells.Item(209, 4).Value = 45141
$ws.Cells.Item(209, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(209, 5).Value = 8
$ws.Cells.Item(209, 6).Value = "Fruta"
$ws.Cells.Item(209, 7).Value = 100108
$ws.Cells.Item(209, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(209, 9).Value = 100108005
$ws.Cells.Item(209, 10).Value = "Piña"
$ws.Cells.Item(209, 11).Value = "Caramelo"
$ws.Cells.Item(209, 12).Value = "Tercera"
$ws.Cells.Item(209, 13).Value = 100
$ws.Cells.Item(209, 14).Value = 18000
$ws.Cells.Item(209, 15).Value = 19000
$ws.Cells.Item(209, 16).Value = 18500
$ws.Cells.Item(209, 17).Value = "`$/caja 16 unidades"
$ws.Cells.Item(209, 18).Value = "Ecuador"
$ws.Cells.Item(209, 19).Value = 1156
$ws.Cells.Item(209, 20).Value = 16
